$wb = $excel.ActiveWorkbook

# Row 48 in ALC (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(48, 8).Value = 2410.6667
$ws.Cells.Item(48, 9).Value = 2410.6667
$ws.Cells.Item(48, 11).Value = 7232.000100000001
$ws.Cells.Item(48, 13).Value = -6940.000100000001

# Row 56 in ALC (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(56, 8).Value = 2410.6667
$ws.Cells.Item(56, 9).Value = 2410.6667
$ws.Cells.Item(56, 11).Value = 7232.000100000001
$ws.Cells.Item(56, 13).Value = -6698.000100000001

# Row 70 in ALC (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2530.5
$ws.Cells.Item(70, 9).Value = 1192.75
$ws.Cells.Item(70, 10).Value = 3868.25
$ws.Cells.Item(70, 11).Value = 3578.25
$ws.Cells.Item(70, 12).Value = 11604.75
$ws.Cells.Item(70, 13).Value = -3308.25
$ws.Cells.Item(70, 14).Value = -12144.75

# Row 73 in ALC (hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 2530.5
$ws.Cells.Item(73, 9).Value = 1192.75
$ws.Cells.Item(73, 10).Value = 3868.25
$ws.Cells.Item(73, 11).Value = 3578.25
$ws.Cells.Item(73, 12).Value = 11604.75
$ws.Cells.Item(73, 13).Value = -2642.25
$ws.Cells.Item(73, 14).Value = -13476.75

# Row 132 in ALC (hunk 4)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1858.8
$ws.Cells.Item(132, 9).Value = 1447.32
$ws.Cells.Item(132, 10).Value = 2887.5
$ws.Cells.Item(132, 11).Value = 4341.96
$ws.Cells.Item(132, 12).Value = 8662.5
$ws.Cells.Item(132, 13).Value = -1811.96
$ws.Cells.Item(132, 14).Value = -13722.5

# Row 137 in ALC (hunk 5)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1949
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 1949
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 5847
$ws.Cells.Item(137, 13).ClearContents()
$ws.Cells.Item(137, 14).Value = -10947

# Row 26 in ARM (hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 92496.25
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()

# Row 32 in ARM (hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9873
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 9873
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 9873
$ws.Cells.Item(32, 13).ClearContents()
$ws.Cells.Item(32, 14).Value = -10447

# Row 74 in ARM (hunk 8)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 10703.583
$ws.Cells.Item(74, 10).Value = 36931.668
$ws.Cells.Item(74, 12).Value = 36931.668
$ws.Cells.Item(74, 14).Value = -38679.668

# Row 77 in ARM (hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 10703.583
$ws.Cells.Item(77, 10).Value = 36931.668
$ws.Cells.Item(77, 12).Value = 184658.34
$ws.Cells.Item(77, 14).Value = -193394.34

# Row 132 in ARM (hunk 10)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3566.6667
$ws.Cells.Item(132, 9).Value = 3380.3333
$ws.Cells.Item(132, 11).Value = 10140.9999
$ws.Cells.Item(132, 13).Value = -7610.999899999999

# Row 105 in BSM (hunk 11)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 4512.6
$ws.Cells.Item(105, 9).Value = 1187.8334
$ws.Cells.Item(105, 11).Value = 1187.8334
$ws.Cells.Item(105, 13).Value = 559.1666

# Row 134 in BSM (hunk 12)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 1976.1915
$ws.Cells.Item(134, 9).Value = 1978.7954
$ws.Cells.Item(134, 10).Value = 1938
$ws.Cells.Item(134, 11).Value = 5936.3862
$ws.Cells.Item(134, 12).Value = 5814
$ws.Cells.Item(134, 13).Value = -3401.3862
$ws.Cells.Item(134, 14).Value = -10884

# Row 3 in CRP (hunk 13)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 853583.2
$ws.Cells.Item(3, 9).Value = 1252874.8
$ws.Cells.Item(3, 10).Value = 55000
$ws.Cells.Item(3, 11).Value = 1252874.8
$ws.Cells.Item(3, 12).Value = 55000
$ws.Cells.Item(3, 13).Value = -1252761.8
$ws.Cells.Item(3, 14).Value = -55226

# Row 31 in CRP (hunk 14)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 23782.936
$ws.Cells.Item(31, 10).Value = 3668.7058
$ws.Cells.Item(31, 12).Value = 3668.7058
$ws.Cells.Item(31, 14).Value = -4258.7058

# Row 34 in CRP (hunk 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 23782.936
$ws.Cells.Item(34, 10).Value = 3668.7058
$ws.Cells.Item(34, 12).Value = 3668.7058
$ws.Cells.Item(34, 14).Value = -4072.7058

# Row 58 in CRP (hunk 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2700.6382
$ws.Cells.Item(58, 9).Value = 2880.1035
$ws.Cells.Item(58, 10).Value = 2411.5
$ws.Cells.Item(58, 11).Value = 2880.1035
$ws.Cells.Item(58, 12).Value = 2411.5
$ws.Cells.Item(58, 13).Value = -2677.1035
$ws.Cells.Item(58, 14).Value = -2817.5

# Row 86 in CRP (hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 5249.5
$ws.Cells.Item(86, 9).Value = 5249.5
$ws.Cells.Item(86, 11).Value = 5249.5
$ws.Cells.Item(86, 13).Value = -4126.5

# Row 89 in CRP (hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 5249.5
$ws.Cells.Item(89, 9).Value = 5249.5
$ws.Cells.Item(89, 11).Value = 26247.5
$ws.Cells.Item(89, 13).Value = -20631.5

# Row 103 in CRP (hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(103, 8).Value = 10000
$ws.Cells.Item(103, 9).Value = 10000
$ws.Cells.Item(103, 10).Value = 10000
$ws.Cells.Item(103, 11).Value = 10000
$ws.Cells.Item(103, 12).Value = 10000
$ws.Cells.Item(103, 13).Value = -8828
$ws.Cells.Item(103, 14).Value = -12344

# Row 121 in CRP (hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(121, 8).Value = 45846.777
$ws.Cells.Item(121, 9).Value = 37296
$ws.Cells.Item(121, 10).Value = 46915.625
$ws.Cells.Item(121, 11).Value = 37296
$ws.Cells.Item(121, 12).Value = 46915.625
$ws.Cells.Item(121, 13).Value = -35986
$ws.Cells.Item(121, 14).Value = -49535.625

# Row 132 in CRP (hunk 21)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 7124.1934
$ws.Cells.Item(132, 9).Value = 4693.875
$ws.Cells.Item(132, 11).Value = 14081.625
$ws.Cells.Item(132, 13).Value = -11551.625

# Row 134 in CRP (hunk 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 19563.96
$ws.Cells.Item(134, 9).Value = 10913
$ws.Cells.Item(134, 11).Value = 32739
$ws.Cells.Item(134, 13).Value = -30204

# Row 136 in CRP (hunk 23)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2700.6382
$ws.Cells.Item(136, 9).Value = 2880.1035
$ws.Cells.Item(136, 10).Value = 2411.5
$ws.Cells.Item(136, 11).Value = 8640.3105
$ws.Cells.Item(136, 12).Value = 7234.5
$ws.Cells.Item(136, 13).Value = -6090.3105
$ws.Cells.Item(136, 14).Value = -12334.5

# Row 12 in CUL (hunk 24)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 433.75
$ws.Cells.Item(12, 9).Value = 145.16667
$ws.Cells.Item(12, 11).Value = 435.50001
$ws.Cells.Item(12, 13).Value = -262.50001

# Row 23 in CUL (hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 1217.625
$ws.Cells.Item(23, 9).Value = 1684.25
$ws.Cells.Item(23, 11).Value = 5052.75
$ws.Cells.Item(23, 13).Value = -4817.75

# Row 69 in CUL (hunk 26)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 4250
$ws.Cells.Item(69, 10).Value = 4500
$ws.Cells.Item(69, 12).Value = 13500
$ws.Cells.Item(69, 14).Value = -15122

# Row 72 in CUL (hunk 27)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 4250
$ws.Cells.Item(72, 10).Value = 4500
$ws.Cells.Item(72, 12).Value = 40500
$ws.Cells.Item(72, 14).Value = -48612

# Row 94 in CUL (hunk 28)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 5880.125
$ws.Cells.Item(94, 9).Value = 3306.4285
$ws.Cells.Item(94, 10).Value = 7881.8887
$ws.Cells.Item(94, 11).Value = 9919.2855
$ws.Cells.Item(94, 12).Value = 23645.6661
$ws.Cells.Item(94, 13).Value = -9243.2855
$ws.Cells.Item(94, 14).Value = -24997.6661

# Row 57 in GSM (hunk 29)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 10672.333
$ws.Cells.Item(57, 10).Value = 14999
$ws.Cells.Item(57, 12).Value = 14999
$ws.Cells.Item(57, 14).Value = -16639

# Row 70 in GSM (hunk 30)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7947.593
$ws.Cells.Item(70, 9).Value = 5231.1304
$ws.Cells.Item(70, 10).Value = 23567.25
$ws.Cells.Item(70, 11).Value = 5231.1304
$ws.Cells.Item(70, 12).Value = 23567.25
$ws.Cells.Item(70, 13).Value = -4961.1304
$ws.Cells.Item(70, 14).Value = -24107.25

# Row 73 in GSM (hunk 31)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 7947.593
$ws.Cells.Item(73, 9).Value = 5231.1304
$ws.Cells.Item(73, 10).Value = 23567.25
$ws.Cells.Item(73, 11).Value = 5231.1304
$ws.Cells.Item(73, 12).Value = 23567.25
$ws.Cells.Item(73, 13).Value = -4295.1304
$ws.Cells.Item(73, 14).Value = -25439.25

# Row 80 in GSM (hunk 32)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2740.6667
$ws.Cells.Item(80, 9).Value = 2829
$ws.Cells.Item(80, 10).Value = 2299
$ws.Cells.Item(80, 11).Value = 2829
$ws.Cells.Item(80, 12).Value = 2299
$ws.Cells.Item(80, 13).Value = -1831
$ws.Cells.Item(80, 14).Value = -4295

# Row 83 in GSM (hunk 33)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 2740.6667
$ws.Cells.Item(83, 9).Value = 2829
$ws.Cells.Item(83, 10).Value = 2299
$ws.Cells.Item(83, 11).Value = 14145
$ws.Cells.Item(83, 12).Value = 11495
$ws.Cells.Item(83, 13).Value = -9153
$ws.Cells.Item(83, 14).Value = -21479

# Row 132 in GSM (hunk 34)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4134.476
$ws.Cells.Item(132, 9).Value = 4132.1055
$ws.Cells.Item(132, 10).Value = 4157
$ws.Cells.Item(132, 11).Value = 12396.3165
$ws.Cells.Item(132, 12).Value = 12471
$ws.Cells.Item(132, 13).Value = -9866.316499999999
$ws.Cells.Item(132, 14).Value = -17531

# Row 132 in LTW (hunk 35)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 6285.4287
$ws.Cells.Item(132, 9).Value = 5999.8
$ws.Cells.Item(132, 10).Value = 6999.5
$ws.Cells.Item(132, 11).Value = 17999.4
$ws.Cells.Item(132, 12).Value = 20998.5
$ws.Cells.Item(132, 13).Value = -15469.4
$ws.Cells.Item(132, 14).Value = -26058.5

# Row 136 in LTW (hunk 36)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 6216.727
$ws.Cells.Item(136, 9).Value = 5935.5
$ws.Cells.Item(136, 11).Value = 17806.5

# Row 113 in WVR (hunk 37)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1130.16
$ws.Cells.Item(113, 9).Value = 763.4666999999999
$ws.Cells.Item(113, 10).Value = 1680.2
$ws.Cells.Item(113, 11).Value = 2290.4001
$ws.Cells.Item(113, 12).Value = 5040.6
$ws.Cells.Item(113, 13).Value = -120.4000999999998
$ws.Cells.Item(113, 14).Value = -9380.6

# Row 132 in WVR (hunk 38)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1597.4445
$ws.Cells.Item(132, 9).Value = 1544.3529
$ws.Cells.Item(132, 11).Value = 4633.0587
$ws.Cells.Item(132, 13).Value = -2103.0587

# Row 136 in WVR (hunk 39)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1824.125
$ws.Cells.Item(136, 9).Value = 1475.3438
$ws.Cells.Item(136, 11).Value = 4426.0314
$ws.Cells.Item(136, 13).Value = -1876.0314
